{"js": "// Part 1: merge the \"  12  \" / \"Distance_Range_4_point_corre_function_average\"\n// runs into a single run, keeping the first run's formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst mergedText = \"  12  Distance_Range_4_point_corre_function_average\";\nconst target = paragraphs.items.find((p) => p.text === mergedText);\nif (target) {\n  target.getRange().insertText(mergedText, Word.InsertLocation.replace);\n}\n\n// Part 2: append new paragraphs after the \"mfrequency ...\" paragraph\n// (which is the last paragraph in the body), before the section break.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\n  \"total angular momentum J + angular momentum projection in z direction M.\",\n  Word.InsertLocation.end\n);\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nbody.insertParagraph(\"initial vibrational state    \", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Part 1: merge the \"  12  \" / \"Distance_Range_4_point_corre_function_average\"\n# runs into a single run (same formatting carried over from the first run).\n$d = $word.ActiveDocument\n\n$mergedText = \"  12  Distance_Range_4_point_corre_function_average\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $mergedText\n$find.Replacement.Text = $mergedText\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Part 2: append new paragraphs after the \"mfrequency ...\" paragraph\n# (the last paragraph in the body), before the section break.\nfunction Add-ParagraphAtEnd($text) {\n    $endRange = $d.Range()\n    $endRange.Collapse(0)  # wdCollapseEnd\n    $endRange.InsertParagraphAfter()\n    if ($text) {\n        $lastPara = $d.Paragraphs.Last\n        $lastPara.Range.Text = $text\n    }\n}\n\nAdd-ParagraphAtEnd \"\"\nAdd-ParagraphAtEnd \"total angular momentum J + angular momentum projection in z direction M.\"\nAdd-ParagraphAtEnd \"\"\nAdd-ParagraphAtEnd \"initial vibrational state    \"\n"}
